# Hortaliza, Vega Modelo de Temuco - Puerro
# Weekly price-data refresh: a new daily record is inserted at row 163
# (pushing the existing rows 163-231 down to 164-232) and the dimension
# grows from A1:R231 to A1:R232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 163, shifting rows 163:231 down to 164:232.
$ws.Rows(163).Insert()

# Populate the newly inserted row 163 with the new record's data.
$ws.Cells.Item(163, 1).Value  = 10
$ws.Cells.Item(163, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(163, 3).Value  = "La Araucanía"
$ws.Cells.Item(163, 4).Value  = 44795
$ws.Cells.Item(163, 5).Value  = 9
$ws.Cells.Item(163, 6).Value  = 100112005
$ws.Cells.Item(163, 7).Value  = "Puerro"
$ws.Cells.Item(163, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(163, 9).Value  = "Primera"
$ws.Cells.Item(163, 10).Value = 90
$ws.Cells.Item(163, 11).Value = 17000
$ws.Cells.Item(163, 12).Value = 18000
$ws.Cells.Item(163, 13).Value = 17444
$ws.Cells.Item(163, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(163, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(163, 16).Value = 1454
$ws.Cells.Item(163, 17).Value = 12
$ws.Cells.Item(163, 18).Value = "Hortaliza"
